$wb = $excel.ActiveWorkbook

# Sheet "存款" (deposits) - remove stray internal space in two institution names
$wsDeposit = $wb.Worksheets.Item("存款")
$wsDeposit.Range("B2").Value = "中國信託商業銀行斗六分行"
$wsDeposit.Range("B10").Value = "中華郵政股份有限公司斗六西平郵局政治獻金專戶"

# Sheet "保險" (insurance) - remove stray internal space in remark text
$wsInsurance = $wb.Worksheets.Item("保險")
$wsInsurance.Range("E2").Value = "保險年齡至保險人四If四歲繳費期間六年年繳二十萬"

# Sheet "債務" (debt) - remove stray internal spaces in creditor address and date
$wsDebt = $wb.Worksheets.Item("債務")
$wsDebt.Range("D2").Value = "斗六市農會雲林縣斗六市民生路"
$wsDebt.Range("F2").Value = "99年12月08日"
